$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The website text in column C for "Drone Volt" rows carries trailing blank
# lines, matching the shared string already used by rows 6 and 7.
$nl = "`n"
$websiteText = "https://www.dronevolt.com/" + $nl + $nl + $nl + $nl + $nl + $nl + $nl + $nl

# Add a new data row (row 8) for contact "Jimmy" at Drone Volt.
$ws.Range("A8").Value = "yes"
$ws.Range("B8").Value = "Drone Volt "
$ws.Range("C8").Value = $websiteText
$ws.Range("D8").Value = "jimm123@dronevolt.com"
$ws.Range("E8").Value = "Jimmy"

# Hyperlink the website and e-mail cells, matching the style used by the other rows.
$ws.Hyperlinks.Add($ws.Range("C8"), "https://www.dronevolt.com/") | Out-Null
$ws.Range("C8").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:jimm123@dronevolt.com") | Out-Null
$ws.Range("D8").Style = "Hyperlink"

# The multi-line website text otherwise forces an explicit row height; restore
# the row to the sheet's standard auto height, same as the rest of the rows.
$ws.Rows.Item(8).AutoFit() | Out-Null

# Move the active selection from E14 to E13.
$ws.Range("E13").Select() | Out-Null
